$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(93, 2).Value = 6236612
$ws.Cells.Item(93, 6).Value = "Zamora"
$ws.Cells.Item(93, 7).Value = "Carabobo"
$ws.Cells.Item(93, 9).Value = 2
$ws.Cells.Item(93, 10).Value = "A"
$ws.Cells.Item(93, 11).Value = 3.2
$ws.Cells.Item(93, 13).Value = 2.15
$ws.Cells.Item(93, 14).Value = 4.5
$ws.Cells.Item(93, 15).Value = 3.3
$ws.Cells.Item(93, 16).Value = 1.75
$ws.Cells.Item(93, 17).Value = 0.5
$ws.Cells.Item(93, 18).Value = 2
$ws.Cells.Item(93, 19).Value = 1.8
$ws.Cells.Item(93, 21).Value = 1.925
$ws.Cells.Item(93, 22).Value = 1.875
$ws.Cells.Item(93, 24).Value = -1
$ws.Cells.Item(93, 25).Value = 0.75
$ws.Cells.Item(93, 26).Value = -1
$ws.Cells.Item(93, 27).Value = 0.8
$ws.Cells.Item(93, 28).Value = -0.5
$ws.Cells.Item(93, 29).Value = 0.4375
$ws.Cells.Item(94, 2).Value = 6236611
$ws.Cells.Item(94, 6).Value = "Mineros"
$ws.Cells.Item(94, 7).Value = "Monagas"
$ws.Cells.Item(94, 8).Value = 2
$ws.Cells.Item(94, 9).Value = 1
$ws.Cells.Item(94, 11).Value = 3.2
$ws.Cells.Item(94, 12).Value = 3.4
$ws.Cells.Item(94, 13).Value = 2
$ws.Cells.Item(94, 14).Value = 4.2
$ws.Cells.Item(94, 15).Value = 3.8
$ws.Cells.Item(94, 16).Value = 1.65
$ws.Cells.Item(94, 17).Value = 0.75
$ws.Cells.Item(94, 18).Value = 1.95
$ws.Cells.Item(94, 19).Value = 1.85
$ws.Cells.Item(94, 21).Value = 1.825
$ws.Cells.Item(94, 22).Value = 1.975
$ws.Cells.Item(94, 23).Value = 3.2
$ws.Cells.Item(94, 26).Value = 0.95
$ws.Cells.Item(94, 27).Value = -1
$ws.Cells.Item(94, 28).Value = 0.825
$ws.Cells.Item(94, 29).Value = -1
$ws.Cells.Item(95, 2).Value = 6236254
$ws.Cells.Item(95, 6).Value = "Academia Puerto Cabello"
$ws.Cells.Item(95, 7).Value = "Estudiantes Merida"
$ws.Cells.Item(95, 9).Value = 0
$ws.Cells.Item(95, 10).Value = "H"
$ws.Cells.Item(95, 11).Value = 1.727
$ws.Cells.Item(95, 12).Value = 3.4
$ws.Cells.Item(95, 13).Value = 4.333
$ws.Cells.Item(95, 14).Value = 1.666
$ws.Cells.Item(95, 15).Value = 3.4
$ws.Cells.Item(95, 16).Value = 4.75
$ws.Cells.Item(95, 17).Value = -0.75
$ws.Cells.Item(95, 18).Value = 1.875
$ws.Cells.Item(95, 19).Value = 1.925
$ws.Cells.Item(95, 21).Value = 1.9
$ws.Cells.Item(95, 22).Value = 1.9
$ws.Cells.Item(95, 23).Value = 0.6659999999999999
$ws.Cells.Item(95, 25).Value = -1
$ws.Cells.Item(95, 26).Value = 0.4375
$ws.Cells.Item(95, 27).Value = -0.5
$ws.Cells.Item(95, 28).Value = -1
$ws.Cells.Item(95, 29).Value = 0.8999999999999999
$ws.Cells.Item(96, 2).Value = 6236251
$ws.Cells.Item(96, 6).Value = "Angostura FC"
$ws.Cells.Item(96, 7).Value = "Portuguesa"
$ws.Cells.Item(96, 8).Value = 1
$ws.Cells.Item(96, 9).Value = 2
$ws.Cells.Item(96, 10).Value = "A"
$ws.Cells.Item(96, 11).Value = 3.1
$ws.Cells.Item(96, 12).Value = 3.2
$ws.Cells.Item(96, 13).Value = 2.15
$ws.Cells.Item(96, 14).Value = 4
$ws.Cells.Item(96, 15).Value = 3.6
$ws.Cells.Item(96, 16).Value = 1.75
$ws.Cells.Item(96, 18).Value = 1.8
$ws.Cells.Item(96, 19).Value = 2
$ws.Cells.Item(96, 21).Value = 1.95
$ws.Cells.Item(96, 22).Value = 1.85
$ws.Cells.Item(96, 23).Value = -1
$ws.Cells.Item(96, 25).Value = 0.75
$ws.Cells.Item(96, 26).Value = -0.5
$ws.Cells.Item(96, 27).Value = 0.5
$ws.Cells.Item(96, 28).Value = 0.95
$ws.Cells.Item(97, 2).Value = 6236255
$ws.Cells.Item(97, 6).Value = "Deportivo Rayo Zuliano"
$ws.Cells.Item(97, 7).Value = "Caracas"
$ws.Cells.Item(97, 11).Value = 3.75
$ws.Cells.Item(97, 12).Value = 3.1
$ws.Cells.Item(97, 13).Value = 1.95
$ws.Cells.Item(97, 14).Value = 2.9
$ws.Cells.Item(97, 15).Value = 2.875
$ws.Cells.Item(97, 16).Value = 2.45
$ws.Cells.Item(97, 17).Value = 0.25
$ws.Cells.Item(97, 21).Value = 1.85
$ws.Cells.Item(97, 22).Value = 1.95
$ws.Cells.Item(97, 24).Value = 1.875
$ws.Cells.Item(97, 26).Value = 0.3875
$ws.Cells.Item(97, 27).Value = -0.5
$ws.Cells.Item(97, 29).Value = 0.95
$ws.Cells.Item(98, 2).Value = 6236252
$ws.Cells.Item(98, 6).Value = "Deportivo Tachira"
$ws.Cells.Item(98, 7).Value = "CD Hermanos Colmenares"
$ws.Cells.Item(98, 11).Value = 1.363
$ws.Cells.Item(98, 12).Value = 4.2
$ws.Cells.Item(98, 13).Value = 7.5
$ws.Cells.Item(98, 14).Value = 1.333
$ws.Cells.Item(98, 15).Value = 4.5
$ws.Cells.Item(98, 16).Value = 8
$ws.Cells.Item(98, 17).Value = -1.5
$ws.Cells.Item(98, 18).Value = 2
$ws.Cells.Item(98, 19).Value = 1.8
$ws.Cells.Item(98, 21).Value = 1.925
$ws.Cells.Item(98, 22).Value = 1.875
$ws.Cells.Item(98, 23).Value = 0.333
$ws.Cells.Item(98, 26).Value = -1
$ws.Cells.Item(98, 27).Value = 0.8
$ws.Cells.Item(98, 29).Value = 0.875
$ws.Cells.Item(99, 2).Value = 6236253
$ws.Cells.Item(99, 6).Value = "Deportivo La Guaira"
$ws.Cells.Item(99, 7).Value = "UCV"
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 10).Value = "D"
$ws.Cells.Item(99, 11).Value = 1.833
$ws.Cells.Item(99, 12).Value = 3.25
$ws.Cells.Item(99, 13).Value = 4
$ws.Cells.Item(99, 14).Value = 2
$ws.Cells.Item(99, 15).Value = 3.2
$ws.Cells.Item(99, 16).Value = 3.5
$ws.Cells.Item(99, 17).Value = -0.25
$ws.Cells.Item(99, 18).Value = 1.775
$ws.Cells.Item(99, 19).Value = 2.025
$ws.Cells.Item(99, 21).Value = 1.9
$ws.Cells.Item(99, 22).Value = 1.9
$ws.Cells.Item(99, 24).Value = 2.2
$ws.Cells.Item(99, 25).Value = -1
$ws.Cells.Item(99, 26).Value = -0.5
$ws.Cells.Item(99, 27).Value = 0.5125
$ws.Cells.Item(99, 28).Value = -1
$ws.Cells.Item(99, 29).Value = 0.8999999999999999
$ws.Cells.Item(100, 2).Value = 6236614
$ws.Cells.Item(100, 6).Value = "Mineros"
$ws.Cells.Item(100, 7).Value = "Angostura FC"
$ws.Cells.Item(100, 8).Value = 1
$ws.Cells.Item(100, 11).Value = 2.45
$ws.Cells.Item(100, 12).Value = 3.3
$ws.Cells.Item(100, 13).Value = 2.55
$ws.Cells.Item(100, 14).Value = 1.8
$ws.Cells.Item(100, 15).Value = 3.75
$ws.Cells.Item(100, 16).Value = 3.6
$ws.Cells.Item(100, 17).Value = -0.5
$ws.Cells.Item(100, 18).Value = 1.825
$ws.Cells.Item(100, 19).Value = 1.975
$ws.Cells.Item(100, 21).Value = 1.8
$ws.Cells.Item(100, 22).Value = 2
$ws.Cells.Item(100, 25).Value = 2.6
$ws.Cells.Item(100, 27).Value = 0.9750000000000001
$ws.Cells.Item(100, 28).Value = 0.4
$ws.Cells.Item(100, 29).Value = -0.5
$ws.Cells.Item(101, 2).Value = 6236257
$ws.Cells.Item(101, 6).Value = "CD Hermanos Colmenares"
$ws.Cells.Item(101, 7).Value = "Zamora"
$ws.Cells.Item(101, 8).Value = 0
$ws.Cells.Item(101, 11).Value = 2.3
$ws.Cells.Item(101, 12).Value = 3.2
$ws.Cells.Item(101, 13).Value = 2.8
$ws.Cells.Item(101, 14).Value = 1.666
$ws.Cells.Item(101, 15).Value = 3.8
$ws.Cells.Item(101, 16).Value = 4.2
$ws.Cells.Item(101, 17).Value = -0.75
$ws.Cells.Item(101, 18).Value = 1.9
$ws.Cells.Item(101, 19).Value = 1.9
$ws.Cells.Item(101, 21).Value = 1.9
$ws.Cells.Item(101, 22).Value = 1.9
$ws.Cells.Item(101, 25).Value = 3.2
$ws.Cells.Item(101, 27).Value = 0.8999999999999999
$ws.Cells.Item(101, 28).Value = -1
$ws.Cells.Item(101, 29).Value = 0.8999999999999999
$ws.Cells.Item(162, 2).Value = 7952893
$ws.Cells.Item(162, 6).Value = "UCV"
$ws.Cells.Item(162, 7).Value = "Deportivo La Guaira"
$ws.Cells.Item(162, 8).Value = 1
$ws.Cells.Item(162, 9).Value = 1
$ws.Cells.Item(162, 10).Value = "D"
$ws.Cells.Item(162, 11).Value = 2.1
$ws.Cells.Item(162, 12).Value = 3
$ws.Cells.Item(162, 13).Value = 3.25
$ws.Cells.Item(162, 14).Value = 2.25
$ws.Cells.Item(162, 15).Value = 3.1
$ws.Cells.Item(162, 16).Value = 2.9
$ws.Cells.Item(162, 17).Value = -0.25
$ws.Cells.Item(162, 18).Value = 2.025
$ws.Cells.Item(162, 19).Value = 1.775
$ws.Cells.Item(162, 21).Value = 1.8
$ws.Cells.Item(162, 22).Value = 2
$ws.Cells.Item(162, 23).Value = -1
$ws.Cells.Item(162, 24).Value = 2.1
$ws.Cells.Item(162, 26).Value = -0.5
$ws.Cells.Item(162, 27).Value = 0.3875
$ws.Cells.Item(163, 2).Value = 7952905
$ws.Cells.Item(163, 6).Value = "Angostura FC"
$ws.Cells.Item(163, 7).Value = "Deportivo Tachira"
$ws.Cells.Item(163, 8).Value = 2
$ws.Cells.Item(163, 9).Value = 0
$ws.Cells.Item(163, 10).Value = "H"
$ws.Cells.Item(163, 11).Value = 3.6
$ws.Cells.Item(163, 12).Value = 3.6
$ws.Cells.Item(163, 13).Value = 1.8
$ws.Cells.Item(163, 14).Value = 3.75
$ws.Cells.Item(163, 15).Value = 2.875
$ws.Cells.Item(163, 16).Value = 2.1
$ws.Cells.Item(163, 17).Value = 0.25
$ws.Cells.Item(163, 18).Value = 1.95
$ws.Cells.Item(163, 19).Value = 1.85
$ws.Cells.Item(163, 21).Value = 2.025
$ws.Cells.Item(163, 22).Value = 1.775
$ws.Cells.Item(163, 23).Value = 2.75
$ws.Cells.Item(163, 24).Value = -1
$ws.Cells.Item(163, 26).Value = 0.95
$ws.Cells.Item(163, 27).Value = -1
$ws.Cells.Item(196, 2).Value = 7977874
$ws.Cells.Item(196, 5).Value = 45394.83333333334
$ws.Cells.Item(196, 6).Value = "Monagas"
$ws.Cells.Item(196, 7).Value = "Zamora"
$ws.Cells.Item(196, 8).Value = 1
$ws.Cells.Item(196, 9).Value = 0
$ws.Cells.Item(196, 10).Value = "H"
$ws.Cells.Item(196, 13).Value = 4.75
$ws.Cells.Item(196, 14).Value = 1.615
$ws.Cells.Item(196, 15).Value = 3.4
$ws.Cells.Item(196, 16).Value = 5
$ws.Cells.Item(196, 17).Value = -0.75
$ws.Cells.Item(196, 18).Value = 1.8
$ws.Cells.Item(196, 19).Value = 2
$ws.Cells.Item(196, 20).Value = 2.5
$ws.Cells.Item(196, 21).Value = 1.975
$ws.Cells.Item(196, 22).Value = 1.825
$ws.Cells.Item(196, 23).Value = 0.615
$ws.Cells.Item(196, 24).Value = -1
$ws.Cells.Item(196, 25).Value = -1
$ws.Cells.Item(196, 26).Value = 0.4
$ws.Cells.Item(196, 27).Value = -0.5
$ws.Cells.Item(196, 28).Value = -1
$ws.Cells.Item(196, 29).Value = 0.825
$ws.Cells.Item(197, 2).Value = 7977873
$ws.Cells.Item(197, 5).Value = 45396.78125
$ws.Cells.Item(197, 6).Value = "Deportivo Tachira"
$ws.Cells.Item(197, 7).Value = "Carabobo"
$ws.Cells.Item(197, 11).Value = 1.909
$ws.Cells.Item(197, 12).Value = 3.2
$ws.Cells.Item(197, 13).Value = 3.75
$ws.Cells.Item(197, 14).Value = 1.95
$ws.Cells.Item(197, 15).Value = 2.9
$ws.Cells.Item(197, 16).Value = 3.8
$ws.Cells.Item(197, 17).Value = -0.5
$ws.Cells.Item(197, 18).Value = 2.025
$ws.Cells.Item(197, 19).Value = 1.775
$ws.Cells.Item(197, 20).Value = 2
$ws.Cells.Item(197, 21).Value = 1.975
$ws.Cells.Item(197, 22).Value = 1.825
$ws.Cells.Item(198, 2).Value = 7977872
$ws.Cells.Item(198, 5).Value = 45396.89583333334
$ws.Cells.Item(198, 6).Value = "Caracas"
$ws.Cells.Item(198, 7).Value = "UCV"
$ws.Cells.Item(198, 11).Value = 2.45
$ws.Cells.Item(198, 12).Value = 2.875
$ws.Cells.Item(198, 13).Value = 2.875
$ws.Cells.Item(198, 14).Value = 2.15
$ws.Cells.Item(198, 15).Value = 2.9
$ws.Cells.Item(198, 16).Value = 3.5
$ws.Cells.Item(198, 17).Value = -0.25
$ws.Cells.Item(198, 18).Value = 1.85
$ws.Cells.Item(198, 19).Value = 1.95
$ws.Cells.Item(198, 20).Value = 2
$ws.Cells.Item(198, 21).Value = 1.925
$ws.Cells.Item(198, 22).Value = 1.875
$ws.Cells.Item(199, 2).Value = 7977383
$ws.Cells.Item(199, 5).Value = 45397.83333333334
$ws.Cells.Item(199, 6).Value = "Academia Puerto Cabello"
$ws.Cells.Item(199, 7).Value = "Metropolitanos FC"
$ws.Cells.Item(199, 11).Value = 1.666
$ws.Cells.Item(199, 12).Value = 3.6
$ws.Cells.Item(199, 13).Value = 4.5
$ws.Cells.Item(199, 14).Value = 1.7
$ws.Cells.Item(199, 15).Value = 3.6
$ws.Cells.Item(199, 16).Value = 4.333
$ws.Cells.Item(199, 17).Value = -0.75
$ws.Cells.Item(199, 18).Value = 1.95
$ws.Cells.Item(199, 19).Value = 1.85
$ws.Cells.Item(199, 20).Value = 2.25

$ws.Rows(200).Delete()

